$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Row 7: "Experimental" property gets value "true"
# (Set via apostrophe-prefix so it stays text "true" instead of being
# auto-coerced to a boolean by Excel, then repair the cell's number
# format/style back from the quote-prefix style to the plain one used by
# the rest of the column.)
$ws.Range("B7").Value = "'true"
$ws.Range("B6").Copy()
$ws.Range("B7").PasteSpecial(-4122)  # xlPasteFormats

# Row 8: "Date" property value is refreshed to the new timestamp
$ws.Range("B8").Value = "2025-07-21T12:46:15+00:00"

# Row 18: "Compositional" property gets value "false"
$ws.Range("B18").Value = "'false"
$ws.Range("B17").Copy()
$ws.Range("B18").PasteSpecial(-4122)  # xlPasteFormats

$excel.CutCopyMode = $false
